$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update Price (D) and Volume/1h change (E) columns with refreshed crypto market data.
# Column D values are forced to Text format before assignment so that numeric-looking
# price strings (e.g. "547.10", "1.00") are preserved exactly instead of being
# auto-converted to numbers (which would drop trailing zeros / change formatting).
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "60.386.28"
$ws.Range("E2").Value = "  +4.29%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.347.58"
$ws.Range("E3").Value = "  +2.69%  "
$ws.Range("E4").Value = "  +0.05%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "547.10"
$ws.Range("E5").Value = "  +3.22%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.53"
$ws.Range("E6").Value = "  +0.94%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "1.00"
$ws.Range("E7").Value = "  +0.00%  "
$ws.Range("E8").Value = "  +0.35%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.345.25"
$ws.Range("E9").Value = "  +2.65%  "
$ws.Range("E10").Value = "  +2.30%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "5.53"
$ws.Range("E11").Value = "  +1.29%  "
$ws.Range("E13").Value = "  +2.10%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "23.95"
$ws.Range("E14").Value = "  +2.31%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.763.70"
$ws.Range("E15").Value = "  +2.65%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "60.319.39"
$ws.Range("E16").Value = "  +4.35%  "
$ws.Range("E17").Value = "  +1.78%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.349.71"
$ws.Range("E18").Value = "  +2.84%  "
$ws.Range("E19").Value = "  +1.63%  "
$ws.Range("E20").Value = "  -0.24%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "6.83"
$ws.Range("E21").Value = "  +7.03%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "314.93"
$ws.Range("E22").Value = "  +1.46%  "
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "63.63"
$ws.Range("E24").Value = "  +2.12%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.174"
$ws.Range("E25").Value = "  +3.70%  "
$ws.Range("E26").Value = "  -0.23%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "7.90"
$ws.Range("E27").Value = "  -0.80%  "
$ws.Range("E28").Value = "  +8.65%  "
$ws.Range("E29").Value = "  +2.91%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "171.67"
$ws.Range("E30").Value = "  +0.67%  "
$ws.Range("E31").Value = "  +12.69%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.0₃0731"
$ws.Range("E32").Value = "  +2.27%  "
$ws.Range("E33").Value = "  +4.18%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "1.41"
$ws.Range("E34").Value = "  +14.98%  "
$ws.Range("E36").Value = "  +1.81%  "
$ws.Range("E38").Value = "  +0.06%  "
$ws.Range("E39").Value = "  +7.82%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "324.84"
$ws.Range("E40").Value = "  +13.92%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "38.06"
$ws.Range("E41").Value = "  -1.11%  "
$ws.Range("E42").Value = "  +3.22%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "141.56"
$ws.Range("E43").Value = "  +0.41%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "3.47"
$ws.Range("E44").Value = "  +1.82%  "
$ws.Range("E45").Value = "  +0.47%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "19.58"
$ws.Range("E46").Value = "  +9.47%  "
$ws.Range("E47").Value = "  +1.06%  "
$ws.Range("E48").Value = "  +1.95%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0213"
$ws.Range("E49").Value = "  +2.42%  "

# Rows 50 and 51 swapped coin order (WhiteBITCoin now ranks above BabyDogeCoin),
# with refreshed price/volume figures for both coins.
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "11.03"
$ws.Range("E50").Value = "  +1.03%  "

$ws.Range("B51").Value = "BabyDogeCoin"
$ws.Range("C51").Value = "https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0₆0210"
$ws.Range("E51").Value = "  +14.32%  "
